# Fix wrong formula for delta NPV calculation in cell C14:
#   was: =-E9-E3
#   now: =E9-E3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Formula = "=E9-E3"

# Update the active selection to reflect where the author ended up (C13)
$ws.Range("C13").Select()
